$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the old O column (PREMIUM IN USD),
# pushing PREMIUM IN USD / Net in USD from O/P to R/S.
$ws.Range("O1:Q1").EntireColumn.Insert()

# New header cells
$ws.Range("O1").Value = "Country of residence"
$ws.Range("P1").Value = "Deductible"
$ws.Range("Q1").Value = "Sports Activities"

# New data cells for row 2
$ws.Range("N2").Value = "Lebanon"
$ws.Range("O2").Value = "Lebanon"
$ws.Range("P2").Value = "Yes"
$ws.Range("Q2").Value = "Yes"

$ws.Range("N5").Select()
